$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.54%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "46.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.69%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.676"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.09%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.52%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.14%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9907"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.98%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.575"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.18%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1154"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.09%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1940"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.21%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.34%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09996"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04673"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.66%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.68%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001277"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.51%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006039"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.25%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.371"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.45%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.481"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.14%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.13%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1401"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.38%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2649"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04219"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.33%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001309"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.93%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004639"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "7.09%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "10.64%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003745"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.11%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02805"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "9.67%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05784"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.20%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007753"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.85%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1437"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007287"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.37%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.82%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009043"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.21%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3413"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007385"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.97%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.27%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005808"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.06%"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003505"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "12.11%"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003503"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.74%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.27%"
